# Update "想去人数" (number of people interested) values across sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 163
$ws1.Range("F6").Value = 1366
$ws1.Range("F10").Value = 470
$ws1.Range("F26").Value = 470
$ws1.Range("F28").Value = 354

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 56
$ws2.Range("F11").Value = 161

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 163
$ws4.Range("F7").Value = 1366
$ws4.Range("F12").Value = 56
$ws4.Range("F16").Value = 470
$ws4.Range("F35").Value = 161
$ws4.Range("F40").Value = 470
$ws4.Range("F42").Value = 354
